$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '62.337.16'
$ws.Cells.Item(2, 5).Value = '  +2.22%  '
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '2.426.01'
$ws.Cells.Item(3, 5).Value = '  +3.26%  '
$ws.Cells.Item(4, 5).Value = '  -0.04%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '557.33'
$ws.Cells.Item(5, 5).Value = '  +2.17%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '143.61'
$ws.Cells.Item(6, 5).Value = '  +4.55%  '
$ws.Cells.Item(7, 5).Value = '  -0.07%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.535'
$ws.Cells.Item(8, 5).Value = '  +1.64%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '2.425.76'
$ws.Cells.Item(9, 5).Value = '  +3.83%  '
$ws.Cells.Item(10, 5).Value = '  +5.26%  '
$ws.Cells.Item(11, 5).Value = '  -0.05%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '5.41'
$ws.Cells.Item(12, 5).Value = '  +1.91%  '
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.352'
$ws.Cells.Item(13, 5).Value = '  +2.71%  '
$ws.Cells.Item(14, 5).Value = '  +6.35%  '
$ws.Cells.Item(15, 5).Value = '  +9.56%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '2.865.42'
$ws.Cells.Item(16, 5).Value = '  +3.26%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '62.123.31'
$ws.Cells.Item(17, 5).Value = '  +1.96%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '2.427.32'
$ws.Cells.Item(18, 5).Value = '  +3.50%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '11.11'
$ws.Cells.Item(19, 5).Value = '  +4.52%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '325.17'
$ws.Cells.Item(20, 5).Value = '  +1.57%  '
$ws.Cells.Item(21, 5).Value = '  +1.55%  '
$ws.Cells.Item(22, 5).Value = '  +3.17%  '
$ws.Cells.Item(23, 5).Value = '  +0.21%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '1.78'
$ws.Cells.Item(24, 5).Value = '  +2.84%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '65.14'
$ws.Cells.Item(25, 5).Value = '  +2.84%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '9.15'
$ws.Cells.Item(26, 5).Value = '  +8.58%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '571.41'
$ws.Cells.Item(27, 5).Value = '  +13.70%  '
$ws.Cells.Item(28, 2).Value = 'WrappedeETH'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '2.545.14'
$ws.Cells.Item(28, 5).Value = '  +3.16%  '
$ws.Cells.Item(29, 2).Value = 'PEPE'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '0.0₃0948'
$ws.Cells.Item(29, 5).Value = '  +9.23%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '0.999'
$ws.Cells.Item(30, 5).Value = '  +0.02%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '8.41'
$ws.Cells.Item(31, 5).Value = '  +5.55%  '
$ws.Cells.Item(32, 5).Value = '  +5.62%  '
$ws.Cells.Item(33, 5).Value = '  +1.52%  '
$ws.Cells.Item(34, 5).Value = '  +4.07%  '
$ws.Cells.Item(35, 5).Value = '  +4.66%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '5.73'
$ws.Cells.Item(36, 5).Value = '  +8.83%  '
$ws.Cells.Item(37, 2).Value = 'NEARProtocol'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '4.84'
$ws.Cells.Item(37, 5).Value = '  +4.43%  '
$ws.Cells.Item(38, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.999'
$ws.Cells.Item(38, 5).Value = '  -0.07%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.385'
$ws.Cells.Item(39, 5).Value = '  +1.97%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '1.90'
$ws.Cells.Item(40, 5).Value = '  +3.09%  '
$ws.Cells.Item(41, 5).Value = '  +1.79%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '149.43'
$ws.Cells.Item(42, 5).Value = '  +5.98%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.999'
$ws.Cells.Item(43, 5).Value = '  +0.01%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '41.74'
$ws.Cells.Item(44, 5).Value = '  +2.89%  '
$ws.Cells.Item(45, 5).Value = '  +12.30%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '152.37'
$ws.Cells.Item(46, 5).Value = '  +6.74%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '3.65'
$ws.Cells.Item(47, 5).Value = '  +2.33%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '0.0545'
$ws.Cells.Item(48, 5).Value = '  +5.12%  '
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '20.42'
$ws.Cells.Item(49, 5).Value = '  +6.89%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.590'
$ws.Cells.Item(50, 5).Value = '  +3.82%  '
$ws.Cells.Item(51, 5).Value = '  +3.51%  '
